# Update birth_year (column Q) and age_y (column S) for rows 2 through 43:
# birth_year decreases by 1, age_y increases by 1 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 43; $row++) {
    $qCell = $ws.Cells.Item($row, 17)  # column Q = birth_year
    $sCell = $ws.Cells.Item($row, 19)  # column S = age_y

    $qCell.Value = $qCell.Value() - 1
    $sCell.Value = $sCell.Value() + 1
}
